$wb = $excel.ActiveWorkbook

# --- "Tables" sheet: replace the single "Occupied" column with a row of
# half-hour time-slot headers (00:00 .. 23:00) across columns C..AW, and
# mark every slot in row 2 (Table 1's availability row) as unoccupied.
$ws = $wb.Worksheets.Item("Tables")

for ($i = 0; $i -lt 47; $i++) {
    $h = [math]::Floor($i / 2)
    $m = ($i % 2) * 30
    $label = "{0:D2}:{1:D2}" -f $h, $m
    $col = 3 + $i
    $ws.Cells.Item(1, $col).Value = $label
    $ws.Cells.Item(2, $col).Value = $false
}

# Row 3 held the "Table 2" entry that is being removed entirely.
$ws.Rows.Item(3).Delete()

# --- Drop the now-unused "Table 2" worksheet (and its sheet3.xml part).
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Table 2").Delete()
